$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New applicant record appended as row 51.
# Plain text fields can be assigned directly.
$ws.Range("A51").Value = "Ashirboyeva Shoxsanam Shoikromovna"
$ws.Range("B51").Value = "Maktabgacha talim tashkiloti direktori"
$ws.Range("C51").Value = "AD5815225"
$ws.Range("E51").Value = "Toshkent shahri"
$ws.Range("F51").Value = "Mirzo Ulug" + [char]0x02BB + "bek tumani"

# D51 ("217") and H51 ("03-11-2024") look numeric/date, so a plain .Value
# assignment would silently coerce them into a number / date serial. Build
# the text in a scratch cell via TEXT(), copy it, and paste values-only into
# the destination so the cell keeps General/default styling (no NumberFormat
# property is ever touched on the destination, so no stray style entries are
# created) while still landing as a genuine text cell.
$ws.Range("ZZ1").Formula = '=TEXT(217,"0")'
$ws.Range("ZZ1").Copy()
$ws.Range("D51").PasteSpecial(-4163)

$ws.Range("ZZ1").Formula = '=TEXT(DATE(2024,11,3),"dd-mm-yyyy")'
$ws.Range("ZZ1").Copy()
$ws.Range("H51").PasteSpecial(-4163)

$ws.Range("ZZ1").ClearContents()

# G51 is a genuine number, formatted with the fraction number format.
$ws.Range("G51").Value = 998909794434
$ws.Range("G51").NumberFormat = "# ?/?"

$excel.CutCopyMode = 0

$ws.Application.Goto($ws.Range("H48"), $true)
